# Add a "SearchData"/"Director" column (B) next to the existing
# "SectorName"/"Banking and finance" column (A) on the HomePage sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "SearchData"
$ws.Range("B2").Value = "Director"

# Copy header style (bold font + yellow fill) from A1 to B1
$ws.Range("A1").Copy() | Out-Null
$ws.Range("B1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = $false

# Set column B width to match target (~14.63 chars, closest reachable value)
$ws.Range("B1").ColumnWidth = 13.8

# Update the active cell selection
$ws.Range("B3").Select() | Out-Null
